$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.667069666666666
$ws.Range("H2").Value = 11.001209
$ws.Range("I2").Value = 0.01298011522000835
$ws.Range("J2").Value = 0.01298011522000835
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 534.2947811778785
$ws.Range("R2").Value = 4808.653030600906
$ws.Range("S2").Value = 0.003720020220974534
$ws.Range("T2").Value = 0.003720020220974534
$ws.Range("G3").Value = 3.667069666666666
$ws.Range("H3").Value = 11.001209
$ws.Range("I3").Value = 0.01298011522000835
$ws.Range("J3").Value = 0.01298011522000835
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 619.0002693912858
$ws.Range("R3").Value = 5571.002424521573
$ws.Range("S3").Value = 0.004309781042307523
$ws.Range("T3").Value = 0.004309781042307524
$ws.Range("G4").Value = 3.667069666666666
$ws.Range("H4").Value = 11.001209
$ws.Range("I4").Value = 0.01298011522000835
$ws.Range("J4").Value = 0.01298011522000835
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 469.8475352847751
$ws.Range("R4").Value = 4228.627817562976
$ws.Range("S4").Value = 0.003271307139068825
$ws.Range("T4").Value = 0.003271307139068825
$ws.Range("G5").Value = 3.667069666666666
$ws.Range("H5").Value = 11.001209
$ws.Range("I5").Value = 0.01298011522000835
$ws.Range("J5").Value = 0.01298011522000835
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 241.1504580481703
$ws.Range("R5").Value = 2170.354122433533
$ws.Range("S5").Value = 0.001679006817657471
$ws.Range("T5").Value = 0.001679006817657471
$ws.Range("I6").Value = 0.5954329572989919
$ws.Range("J6").Value = 0.595432957298992
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 24509.54527243074
$ws.Range("R6").Value = 220585.9074518766
$ws.Range("S6").Value = 0.1706473789980341
$ws.Range("T6").Value = 0.1706473789980341
$ws.Range("I7").Value = 0.5954329572989919
$ws.Range("J7").Value = 0.595432957298992
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.1977013014011327
$ws.Range("T7").Value = 0.1977013014011327
$ws.Range("I8").Value = 0.5954329572989919
$ws.Range("J8").Value = 0.595432957298992
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 21553.17596742228
$ws.Range("R8").Value = 193978.5837068004
$ws.Range("S8").Value = 0.1500636975122168
$ws.Range("T8").Value = 0.1500636975122168
$ws.Range("I9").Value = 0.5954329572989919
$ws.Range("J9").Value = 0.595432957298992
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 11062.22309708713
$ws.Range("R9").Value = 99560.00787378412
$ws.Range("S9").Value = 0.0770205793876084
$ws.Range("T9").Value = 0.07702057938760841
$ws.Range("G10").Value = 110.4727123333333
$ws.Range("H10").Value = 331.418137
$ws.Range("I10").Value = 0.3910338949346852
$ws.Range("J10").Value = 0.3910338949346853
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 16095.95645231312
$ws.Range("R10").Value = 144863.6080708181
$ws.Range("S10").Value = 0.1120678801064236
$ws.Range("T10").Value = 0.1120678801064236
$ws.Range("G11").Value = 110.4727123333333
$ws.Range("H11").Value = 331.418137
$ws.Range("I11").Value = 0.3910338949346852
$ws.Range("J11").Value = 0.3910338949346853
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 18647.76099464687
$ws.Range("R11").Value = 167829.8489518218
$ws.Range("S11").Value = 0.1298347848785963
$ws.Range("T11").Value = 0.1298347848785963
$ws.Range("G12").Value = 110.4727123333333
$ws.Range("H12").Value = 331.418137
$ws.Range("I12").Value = 0.3910338949346852
$ws.Range("J12").Value = 0.3910338949346853
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 14154.44382686684
$ws.Range("R12").Value = 127389.9944418016
$ws.Range("S12").Value = 0.09855012458948738
$ws.Range("T12").Value = 0.09855012458948739
$ws.Range("G13").Value = 110.4727123333333
$ws.Range("H13").Value = 331.418137
$ws.Range("I13").Value = 0.3910338949346852
$ws.Range("J13").Value = 0.3910338949346853
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 7264.804763096608
$ws.Range("R13").Value = 65383.24286786946
$ws.Range("S13").Value = 0.05058110536017794
$ws.Range("T13").Value = 0.05058110536017794
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1562396666666667
$ws.Range("H14").Value = 0.468719
$ws.Range("I14").Value = 0.0005530325463144183
$ws.Range("J14").Value = 0.0005530325463144184
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 22.76423577980511
$ws.Range("R14").Value = 204.878122018246
$ws.Range("S14").Value = 0.0001584956851519649
$ws.Range("T14").Value = 0.0001584956851519649
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1562396666666667
$ws.Range("H15").Value = 0.468719
$ws.Range("I15").Value = 0.0005530325463144183
$ws.Range("J15").Value = 0.0005530325463144184
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 26.37320927807245
$ws.Range("R15").Value = 237.358883502652
$ws.Range("S15").Value = 0.0001836231145476229
$ws.Range("T15").Value = 0.0001836231145476229
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1562396666666667
$ws.Range("H16").Value = 0.468719
$ws.Range("I16").Value = 0.0005530325463144183
$ws.Range("J16").Value = 0.0005530325463144184
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 20.01838769640178
$ws.Range("R16").Value = 180.165489267616
$ws.Range("S16").Value = 0.0001393777548374184
$ws.Range("T16").Value = 0.0001393777548374184
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1562396666666667
$ws.Range("H17").Value = 0.468719
$ws.Range("I17").Value = 0.0005530325463144183
$ws.Range("J17").Value = 0.0005530325463144184
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 10.27448906260033
$ws.Range("R17").Value = 92.47040156340299
$ws.Range("S17").Value = 0.00007153599177741212
$ws.Range("T17").Value = 0.00007153599177741213
